# Update marksheet corrected/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: correct answers count changed 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total marks changed 69 -> 115
$ws.Range("B12").Value = 115

# "Total" row: corrected/total marks label changed "64/84" -> "115/140"
$ws.Range("E12").Value = "115/140"
